# Natmi following Dr Hou advice
# Rewrites the LR-pairs result table: expands from 3 data rows to 6 data
# rows (ECs/sCs sending clusters x ECs/FAPs/sCs target clusters) and
# refreshes the computed statistics.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("ECs", "Osm", "Il6st", "ECs", 3, 1, 21.90566233333334, 65.716987, 0.9982377076397728, 0.9982377076397729, 3, 1, 37.42645866666667, 112.279376, 0.2415534622699011, 0.2415534622699011, 819.851365884457, 7378.662292960113, 0.2411277744487564, 0.2411277744487565),
    @("ECs", "Osm", "Il6st", "FAPs", 3, 1, 21.90566233333334, 65.716987, 0.9982377076397728, 0.9982377076397729, 3, 1, 94.96115633333334, 284.883469, 0.6128871635375853, 0.6128871635375853, 2080.187025420878, 18721.6832287879, 0.6118070771716018, 0.6118070771716018),
    @("ECs", "Osm", "Il6st", "sCs", 3, 1, 21.90566233333334, 65.716987, 0.9982377076397728, 0.9982377076397729, 3, 1, 22.553069, 67.659207, 0.1455593741925136, 0.1455593741925136, 494.0399140943676, 4446.359226849308, 0.1453028560194146, 0.1453028560194146),
    @("sCs", "Osm", "Il6st", "ECs", 1, 0.3333333333333333, 0.03867233333333333, 0.116017, 0.001762292360227098, 0.001762292360227098, 3, 1, 37.42645866666667, 112.279376, 0.2415534622699011, 0.2415534622699011, 1.447368485043556, 13.026316365392, 0.0004256878211446512, 0.0004256878211446513),
    @("sCs", "Osm", "Il6st", "FAPs", 1, 0.3333333333333333, 0.03867233333333333, 0.116017, 0.001762292360227098, 0.001762292360227098, 3, 1, 94.96115633333334, 284.883469, 0.6128871635375853, 0.6128871635375853, 3.672369491441444, 33.051325422973, 0.001080086365983543, 0.001080086365983543),
    @("sCs", "Osm", "Il6st", "sCs", 1, 0.3333333333333333, 0.03867233333333333, 0.116017, 0.001762292360227098, 0.001762292360227098, 3, 1, 22.553069, 67.659207, 0.1455593741925136, 0.1455593741925136, 0.8721798020576664, 7.849618218518999, 0.000256518173098904, 0.000256518173098904)
)

$startRow = 2
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $values = $data[$i]
    for ($c = 0; $c -lt $values.Count; $c++) {
        $ws.Cells.Item($row, $c + 1).Value = $values[$c]
    }
}
